$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.986.39"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "2.588.94"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'582.23"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").Value = "'147.86"
$ws.Range("E6").Value = "  +1.46%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.598"
$ws.Range("E8").Value = "  +2.78%  "
$ws.Range("E9").Value = "  +2.86%  "
$ws.Range("D10").Value = "'5.66"
$ws.Range("E10").Value = "  +3.27%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "'0.353"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "'27.27"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "3.054.80"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").Value = "62.875.38"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "'0.0000147"
$ws.Range("E16").Value = "  +3.51%  "
$ws.Range("D17").Value = "2.590.85"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").Value = "'11.33"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").Value = "'343.22"
$ws.Range("E19").Value = "  +2.77%  "
$ws.Range("D20").Value = "'4.39"
$ws.Range("E20").Value = "  +2.13%  "
$ws.Range("D21").Value = "'6.69"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "'5.65"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").Value = "'67.21"
$ws.Range("E24").Value = "  +3.42%  "
$ws.Range("D25").Value = "2.725.26"
$ws.Range("E25").Value = "  +2.57%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").Value = "'1.59"
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "'8.34"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").Value = "'7.83"
$ws.Range("E30").Value = "  +8.46%  "
$ws.Range("D31").Value = "'1.44"
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("D32").Value = "'1.93"
$ws.Range("E32").Value = "  +4.98%  "
$ws.Range("D33").Value = "0.0₃0826"
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("D34").Value = "'465.86"
$ws.Range("E34").Value = "  +17.12%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'175.03"
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.61"
$ws.Range("E36").Value = "  +5.31%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "'0.403"
$ws.Range("E38").Value = "  +1.75%  "
$ws.Range("D39").Value = "'19.14"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("D40").Value = "'4.56"
$ws.Range("E40").Value = "  +5.82%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'1.70"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("D43").Value = "'158.83"
$ws.Range("E43").Value = "  +5.65%  "
$ws.Range("D44").Value = "'3.78"
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("D45").Value = "'0.639"
$ws.Range("E45").Value = "  +6.94%  "
$ws.Range("D46").Value = "'21.05"
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("D47").Value = "'0.0546"
$ws.Range("E47").Value = "  +3.10%  "
$ws.Range("D48").Value = "'0.0970"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").Value = "'0.0238"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "'18.59"
$ws.Range("E50").Value = "  +3.21%  "
$ws.Range("D51").Value = "'1.71"
$ws.Range("E51").Value = "  +1.46%  "
